# Edit script: apply the "May 9th" data update to Sheet1.
# The accelerometer/gyroscope sample buffer scrolled forward by one sample:
# the oldest sample (original row 2) was dropped, all remaining samples
# shifted up by one row, and 10 brand-new samples were appended at the
# bottom. The timestamp column (A) keeps counting 0,100,200,... regardless.
# We simply rewrite the full A2:H31 data block with the final values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0, "struggle", -0.2798952907323835, -0.2011211831122636, -0.2037543449550868, -0.380110889673233, -2.916423320770264, -0.522595226764679),
    @(100, "struggle", 0.07621550559997564, 0.3202240578830243, -1.714594541117549, -0.2550362348556518, -3.585931777954102, 1.043205142021179),
    @(200, "struggle", 1.214714780449869, -0.2680620010942265, -5.130570538342003, 0.2286163121461868, -2.957656860351562, 1.703550457954407),
    @(300, "struggle", 2.21673782914877, -2.121318969875573, -6.365980595350262, -0.0294742472469806, -2.758056640625, 1.507310032844544),
    @(400, "struggle", 1.646873444318771, -0.9547826964408157, -2.861810609698295, 0.2777909636497497, 0.1227839142084121, 0.2632829248905182),
    @(500, "struggle", -1.365792244672773, -0.1524508148431776, -1.155832886695863, 0.5624541640281677, 3.851505279541016, -0.9796751141548156),
    @(600, "struggle", -5.90024196356535, -0.5142916264012455, 2.305023960769181, -0.4262312352657318, 4.981453418731689, -0.9717338681221008),
    @(700, "struggle", -3.965735837817196, -0.6069814283400773, -1.651228621602052, -0.7962628602981567, 4.429689407348633, 0.3517056703567505),
    @(800, "struggle", -5.797070525586607, -1.400326563976706, 1.329120457172396, -0.942259669303894, 2.608700037002563, 0.6546948552131653),
    @(900, "struggle", -1.573704838752743, -1.781807705760002, 1.208315253257751, -0.8762862086296082, 1.198364853858948, 0.3744603991508484),
    @(1000, "struggle", -1.48451391607523, -2.143217463046312, -4.305671475827694, -0.3497203588485718, -0.912632703781128, -0.2125810980796814),
    @(1100, "struggle", -3.191123038530348, -1.04998478293419, -3.277032032608988, 1.418276429176331, -3.182302474975586, -0.7964155673980713),
    @(1200, "struggle", -2.79414886981249, 1.120455805212256, -3.891929678618904, 0.7744244337081909, -4.472449779510498, 0.1556179225444793),
    @(1300, "struggle", -0.3785691261291502, -0.5515228509902954, -5.506646633148193, 0.039248090237379, -4.096309661865234, 0.9094256162643432),
    @(1400, "struggle", 0.511223681271077, -3.341538593173027, -3.347421273589136, 0.1157589629292488, -3.650377988815308, 0.6218608021736145),
    @(1500, "struggle", -1.899449691176408, -2.177598688751465, 0.188652187585828, 0.1786780804395675, -0.7522805333137512, 0.2764165103435516),
    @(1600, "struggle", -4.877644158899785, -0.1642134003341204, 1.343364715576172, -0.0583376325666904, 4.005901336669922, -1.094670534133911),
    @(1700, "struggle", -5.74425581097603, -0.5421471670269949, 1.297248423099513, 0.1313360333442688, 5.505575656890869, -0.8086328506469727),
    @(1800, "struggle", -6.533521950244904, -0.2819428481161567, 0.1814587861299474, -0.5644394755363464, 4.288732051849365, 0.6246097087860107),
    @(1900, "struggle", -5.779280513524995, -3.271673738956447, 4.346601516008363, -0.2765692472457886, 2.626567840576172, 0.6087272167205811),
    @(2000, "struggle", 3.140387788414934, 0.259726375341407, -4.785581156611421, -0.3174972236156463, 1.27487576007843, 0.1149953827261924),
    @(2100, "struggle", -1.632258296012878, 0.6425724923610687, -3.22618693113327, -0.2593123018741607, -0.741895854473114, -0.290313720703125),
    @(2200, "struggle", -3.271841421723368, 0.07577018067240526, -1.077775649726385, 0.6250678896903992, -3.214983701705933, -0.8246681094169617),
    @(2300, "struggle", -3.740465611219407, -0.2502757757902145, -5.230584308505059, 0.2973386645317077, -4.603633403778076, 0.0426078513264656),
    @(2400, "struggle", -1.173786669969556, -1.206141140311958, -5.999948702752588, -0.3110831379890442, -4.220009803771973, 1.416138410568237),
    @(2500, "struggle", -0.527452439069747, -1.933494433760643, -1.995455801486972, -0.2417499274015426, -3.320205211639404, 1.446528911590576),
    @(2600, "struggle", -2.818732134997842, -1.367858927696945, 0.861171409487724, 0.3194825351238251, -0.44180828332901, 0.5250386595726013),
    @(2700, "struggle", -4.338251754641532, -0.358771674335002, 0.811524987220763, 0.0740674138069152, 2.8290696144104, -0.7684684991836548),
    @(2800, "struggle", -5.003720842301852, -0.2583636995404971, 1.431181490421301, 0.3081815242767334, 5.233893394470215, -0.9094256162643432),
    @(2900, "struggle", -6.261017680168152, -1.120874315500259, 4.274679899215698, -0.4476115107536316, 4.42144250869751, 0.3985895812511444)
)

$nRows = $data.Count
$nCols = $data[0].Count

$arr = New-Object 'object[,]' $nRows, $nCols
for ($i = 0; $i -lt $nRows; $i++) {
    for ($j = 0; $j -lt $nCols; $j++) {
        $arr[$i, $j] = $data[$i][$j]
    }
}

$startRow = 2
$endRow = $startRow + $nRows - 1
$rng = $ws.Range("A$startRow`:H$endRow")
$rng.Value2 = $arr
